$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cell I7 currently holds the "ser: 101" blog entry; replace it with the
# new "ser: 104" blog entry (life on an ocean).
$ws.Range("I7").Value = "type: blog" + [char]10 + "width: 2" + [char]10 + "height: 1" + [char]10 + "ser: 104"
